# Auto-generated edit script: applies numeric "want-to-go" count updates
# and location/cover-image text updates for the gh-pages data refresh
# (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 3074
$ws1.Range("F5").Value = 1687
$ws1.Range("F6").Value = 2082
$ws1.Range("F8").Value = 304
$ws1.Range("F9").Value = 889
$ws1.Range("F10").Value = 973
$ws1.Range("F11").Value = 217
$ws1.Range("F12").Value = 438
$ws1.Range("F13").Value = 1143
$ws1.Range("F17").Value = 7460
$ws1.Range("F18").Value = 309
$ws1.Range("F19").Value = 2449
$ws1.Range("F20").Value = 203
$ws1.Range("F21").Value = 213
$ws1.Range("F22").Value = 167
$ws1.Range("F24").Value = 516
$ws1.Range("F26").Value = 1126
$ws1.Range("F27").Value = 959
$ws1.Range("F29").Value = 653
$ws1.Range("F31").Value = 1134
$ws1.Range("F32").Value = 1910
$ws1.Range("F34").Value = 27
$ws1.Range("F35").Value = 159
$ws1.Range("F38").Value = 156
$ws1.Range("F39").Value = 305
$ws1.Range("F41").Value = 201

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D2").Value = "湖墅南路186-1 DMT CLUB"
$ws2.Range("F2").Value = 20
$ws2.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202403/Ob3tSTRq1709386334660.jpeg"
$ws2.Range("F8").Value = 18

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("D3").Value = "湖墅南路186-1 DMT CLUB"
$ws4.Range("F3").Value = 20
$ws4.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202403/Ob3tSTRq1709386334660.jpeg"
$ws4.Range("F7").Value = 3074
$ws4.Range("F8").Value = 1687
$ws4.Range("F9").Value = 2082
$ws4.Range("F11").Value = 304
$ws4.Range("F12").Value = 889
$ws4.Range("F14").Value = 973
$ws4.Range("F15").Value = 217
$ws4.Range("F16").Value = 438
$ws4.Range("F17").Value = 1143
$ws4.Range("F21").Value = 7460
$ws4.Range("F22").Value = 309
$ws4.Range("F23").Value = 2449
$ws4.Range("F25").Value = 203
$ws4.Range("F26").Value = 213
$ws4.Range("F27").Value = 167
$ws4.Range("F29").Value = 516
$ws4.Range("F31").Value = 1126
$ws4.Range("F32").Value = 959
$ws4.Range("F34").Value = 653
$ws4.Range("F36").Value = 1134
$ws4.Range("F37").Value = 1910
$ws4.Range("F39").Value = 27
$ws4.Range("F40").Value = 159
$ws4.Range("F43").Value = 156
$ws4.Range("F44").Value = 305
$ws4.Range("F46").Value = 18
$ws4.Range("F49").Value = 201

